$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.068078666666667
$ws.Range("H2").Value = 6.204236
$ws.Range("I2").Value = 0.04745723096963421
$ws.Range("J2").Value = 0.04745723096963421
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.872120333333332
$ws.Range("N2").Value = 17.616361
$ws.Range("O2").Value = 0.07819433676692768
$ws.Range("P2").Value = 0.07819433676692769
$ws.Range("Q2").Value = 12.14400678946622
$ws.Range("R2").Value = 109.296061105196
$ws.Range("S2").Value = 0.003710886700465447
$ws.Range("T2").Value = 0.003710886700465448

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.068078666666667
$ws.Range("H3").Value = 6.204236
$ws.Range("I3").Value = 0.04745723096963421
$ws.Range("J3").Value = 0.04745723096963421
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.62659933333333
$ws.Range("N3").Value = 178.879798
$ws.Range("O3").Value = 0.7939998031155241
$ws.Range("P3").Value = 0.7939998031155242
$ws.Range("Q3").Value = 123.3124980471475
$ws.Range("R3").Value = 1109.812482424328
$ws.Range("S3").Value = 0.03768103204629752
$ws.Range("T3").Value = 0.03768103204629752

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.068078666666667
$ws.Range("H4").Value = 6.204236
$ws.Range("I4").Value = 0.04745723096963421
$ws.Range("J4").Value = 0.04745723096963421
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.226786
$ws.Range("N4").Value = 0.680358
$ws.Range("O4").Value = 0.003019928041555994
$ws.Range("P4").Value = 0.003019928041555994
$ws.Range("Q4").Value = 0.4690112884986667
$ws.Range("R4").Value = 4.221101596488
$ws.Range("S4").Value = 0.0001433174225797979
$ws.Range("T4").Value = 0.0001433174225797979

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.068078666666667
$ws.Range("H5").Value = 6.204236
$ws.Range("I5").Value = 0.04745723096963421
$ws.Range("J5").Value = 0.04745723096963421
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.370985666666668
$ws.Range("N5").Value = 28.112957
$ws.Range("O5").Value = 0.1247859320759922
$ws.Range("P5").Value = 0.1247859320759922
$ws.Range("Q5").Value = 19.37993554287245
$ws.Range("R5").Value = 174.419419885852
$ws.Range("S5").Value = 0.005921994800291447
$ws.Range("T5").Value = 0.005921994800291447

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 26.94592933333334
$ws.Range("H6").Value = 80.837788
$ws.Range("I6").Value = 0.6183416582138921
$ws.Range("J6").Value = 0.6183416582138921
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.872120333333332
$ws.Range("N6").Value = 17.616361
$ws.Range("O6").Value = 0.07819433676692768
$ws.Range("P6").Value = 0.07819433676692769
$ws.Range("Q6").Value = 158.2297395388298
$ws.Range("R6").Value = 1424.067655849468
$ws.Range("S6").Value = 0.04835081585939757
$ws.Range("T6").Value = 0.04835081585939758

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 26.94592933333334
$ws.Range("H7").Value = 80.837788
$ws.Range("I7").Value = 0.6183416582138921
$ws.Range("J7").Value = 0.6183416582138921
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.62659933333333
$ws.Range("N7").Value = 178.879798
$ws.Range("O7").Value = 0.7939998031155241
$ws.Range("P7").Value = 0.7939998031155242
$ws.Range("Q7").Value = 1606.694132022981
$ws.Range("R7").Value = 14460.24718820682
$ws.Range("S7").Value = 0.490963154879957
$ws.Range("T7").Value = 0.4909631548799571

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.94592933333334
$ws.Range("H8").Value = 80.837788
$ws.Range("I8").Value = 0.6183416582138921
$ws.Range("J8").Value = 0.6183416582138921
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.226786
$ws.Range("N8").Value = 0.680358
$ws.Range("O8").Value = 0.003019928041555994
$ws.Range("P8").Value = 0.003019928041555994
$ws.Range("Q8").Value = 6.110959529789334
$ws.Range("R8").Value = 54.998635768104
$ws.Range("S8").Value = 0.001867347312902365
$ws.Range("T8").Value = 0.001867347312902365

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.94592933333334
$ws.Range("H9").Value = 80.837788
$ws.Range("I9").Value = 0.6183416582138921
$ws.Range("J9").Value = 0.6183416582138921
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.370985666666668
$ws.Range("N9").Value = 28.112957
$ws.Range("O9").Value = 0.1247859320759922
$ws.Range("P9").Value = 0.1247859320759922
$ws.Range("Q9").Value = 252.5099175576796
$ws.Range("R9").Value = 2272.589258019116
$ws.Range("S9").Value = 0.0771603401616351
$ws.Range("T9").Value = 0.0771603401616351

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.195412
$ws.Range("H10").Value = 0.586236
$ws.Range("I10").Value = 0.004484216469959312
$ws.Range("J10").Value = 0.004484216469959312
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.872120333333332
$ws.Range("N10").Value = 17.616361
$ws.Range("O10").Value = 0.07819433676692768
$ws.Range("P10").Value = 0.07819433676692769
$ws.Range("Q10").Value = 1.147482778577333
$ws.Range("R10").Value = 10.327345007196
$ws.Range("S10").Value = 0.0003506403327878021
$ws.Range("T10").Value = 0.0003506403327878021

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.195412
$ws.Range("H11").Value = 0.586236
$ws.Range("I11").Value = 0.004484216469959312
$ws.Range("J11").Value = 0.004484216469959312
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.62659933333333
$ws.Range("N11").Value = 178.879798
$ws.Range("O11").Value = 0.7939998031155241
$ws.Range("P11").Value = 0.7939998031155242
$ws.Range("Q11").Value = 11.65175302892533
$ws.Range("R11").Value = 104.865777260328
$ws.Range("S11").Value = 0.003560466994275084
$ws.Range("T11").Value = 0.003560466994275085

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.195412
$ws.Range("H12").Value = 0.586236
$ws.Range("I12").Value = 0.004484216469959312
$ws.Range("J12").Value = 0.004484216469959312
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.226786
$ws.Range("N12").Value = 0.680358
$ws.Range("O12").Value = 0.003019928041555994
$ws.Range("P12").Value = 0.003019928041555994
$ws.Range("Q12").Value = 0.044316705832
$ws.Range("R12").Value = 0.398850352488
$ws.Range("S12").Value = 0.00001354201106203736
$ws.Range("T12").Value = 0.00001354201106203736

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.195412
$ws.Range("H13").Value = 0.586236
$ws.Range("I13").Value = 0.004484216469959312
$ws.Range("J13").Value = 0.004484216469959312
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.370985666666668
$ws.Range("N13").Value = 28.112957
$ws.Range("O13").Value = 0.1247859320759922
$ws.Range("P13").Value = 0.1247859320759922
$ws.Range("Q13").Value = 1.831203051094667
$ws.Range("R13").Value = 16.480827459852
$ws.Range("S13").Value = 0.0005595671318343881
$ws.Range("T13").Value = 0.0005595671318343881

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.36831566666667
$ws.Range("H14").Value = 43.104947
$ws.Range("I14").Value = 0.3297168943465144
$ws.Range("J14").Value = 0.3297168943465145
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.872120333333332
$ws.Range("N14").Value = 17.616361
$ws.Range("O14").Value = 0.07819433676692768
$ws.Range("P14").Value = 0.07819433676692769
$ws.Range("Q14").Value = 84.37247858198521
$ws.Range("R14").Value = 759.352307237867
$ws.Range("S14").Value = 0.02578199387427686
$ws.Range("T14").Value = 0.02578199387427687

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.36831566666667
$ws.Range("H15").Value = 43.104947
$ws.Range("I15").Value = 0.3297168943465144
$ws.Range("J15").Value = 0.3297168943465145
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.62659933333333
$ws.Range("N15").Value = 178.879798
$ws.Range("O15").Value = 0.7939998031155241
$ws.Range("P15").Value = 0.7939998031155242
$ws.Range("Q15").Value = 856.7338013511895
$ws.Range("R15").Value = 7710.604212160706
$ws.Range("S15").Value = 0.2617951491949945
$ws.Range("T15").Value = 0.2617951491949946

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.36831566666667
$ws.Range("H16").Value = 43.104947
$ws.Range("I16").Value = 0.3297168943465144
$ws.Range("J16").Value = 0.3297168943465145
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.226786
$ws.Range("N16").Value = 0.680358
$ws.Range("O16").Value = 0.003019928041555994
$ws.Range("P16").Value = 0.003019928041555994
$ws.Range("Q16").Value = 3.258532836780667
$ws.Range("R16").Value = 29.326795531026
$ws.Range("S16").Value = 0.0009957212950117938
$ws.Range("T16").Value = 0.000995721295011794

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.36831566666667
$ws.Range("H17").Value = 43.104947
$ws.Range("I17").Value = 0.3297168943465144
$ws.Range("J17").Value = 0.3297168943465145
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.370985666666668
$ws.Range("N17").Value = 28.112957
$ws.Range("O17").Value = 0.1247859320759922
$ws.Range("P17").Value = 0.1247859320759922
$ws.Range("Q17").Value = 134.6452801664755
$ws.Range("R17").Value = 1211.807521498279
$ws.Range("S17").Value = 0.04114402998223123
$ws.Range("T17").Value = 0.04114402998223124
